# "Reverted to version 9" -- the title shape on slide 1 reads
# "Version " + "2" (two separate runs). Update the numeric run's
# text from "2" to "6", leaving the "Version " run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the title placeholder shape (falls back to the first shape
# whose text starts with "Version " if the name ever changes).
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text.StartsWith("Version ")) {
        $titleShape = $shape
        break
    }
}
if ($titleShape -eq $null) {
    $titleShape = $s.Shapes.Item(1)
}

$tr = $titleShape.TextFrame.TextRange

# The text is "Version 2": the literal "Version " run occupies
# characters 1-8, and the trailing number is the single character
# at position 9. Grab just that trailing run and replace its text.
$numStart = "Version ".Length + 1
$numRun = $tr.Characters($numStart, $tr.Length - $numStart + 1)
$numRun.Text = "6"
